$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1204.5
$ws.Range("I38").Value = 1204.5
$ws.Range("K38").Value = 3613.5
$ws.Range("M38").Value = -3241.5
$ws.Range("H61").Value = 610.8333
$ws.Range("I61").Value = 333
$ws.Range("K61").Value = 999
$ws.Range("M61").Value = -827
$ws.Range("H62").Value = 5424.75
$ws.Range("I62").Value = 6149.5
$ws.Range("J62").Value = 4700
$ws.Range("K62").Value = 6149.5
$ws.Range("L62").Value = 4700
$ws.Range("M62").Value = -5525.5
$ws.Range("N62").Value = -5948
$ws.Range("H65").Value = 5424.75
$ws.Range("I65").Value = 6149.5
$ws.Range("J65").Value = 4700
$ws.Range("K65").Value = 30747.5
$ws.Range("L65").Value = 23500
$ws.Range("M65").Value = -27627.5
$ws.Range("N65").Value = -29740
$ws.Range("H70").Value = 4333.3335
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 4333.3335
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 13000.0005
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -13540.0005
$ws.Range("H73").Value = 4333.3335
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 4333.3335
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 13000.0005
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -14872.0005
$ws.Range("H76").Value = 4517.1875
$ws.Range("I76").Value = 4605.5386
$ws.Range("K76").Value = 4605.5386
$ws.Range("M76").Value = -4290.5386
$ws.Range("H79").Value = 4517.1875
$ws.Range("I79").Value = 4605.5386
$ws.Range("K79").Value = 4605.5386
$ws.Range("M79").Value = -3513.5386
$ws.Range("H87").Value = 19999.908
$ws.Range("J87").Value = 19999.908
$ws.Range("L87").Value = 19999.908
$ws.Range("N87").Value = -22495.908
$ws.Range("H90").Value = 19999.908
$ws.Range("J90").Value = 19999.908
$ws.Range("L90").Value = 59999.724
$ws.Range("N90").Value = -72479.724
$ws.Range("H112").Value = 1623.3
$ws.Range("J112").Value = 1626.3265
$ws.Range("L112").Value = 4878.979499999999
$ws.Range("N112").Value = -7094.979499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 58901.43
$ws.Range("I99").Value = 68051.664
$ws.Range("K99").Value = 68051.664
$ws.Range("M99").Value = -66553.664
$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 911.4286
$ws.Range("I16").Value = 296.6
$ws.Range("J16").Value = 2448.5
$ws.Range("K16").Value = 296.6
$ws.Range("L16").Value = 2448.5
$ws.Range("M16").Value = -9.600000000000023
$ws.Range("N16").Value = -3022.5
$ws.Range("H43").Value = 14163.223
$ws.Range("I43").Value = 15000
$ws.Range("J43").Value = 14058.625
$ws.Range("K43").Value = 15000
$ws.Range("L43").Value = 14058.625
$ws.Range("N43").Value = -14426.625
$ws.Range("M43").Value = -14816
$ws.Range("H101").Value = 14163.223
$ws.Range("I101").Value = 15000
$ws.Range("J101").Value = 14058.625
$ws.Range("K101").Value = 15000
$ws.Range("L101").Value = 14058.625
$ws.Range("N101").Value = -20548.625
$ws.Range("M101").Value = -11755
$ws.Range("H109").Value = 24505.385
$ws.Range("J109").Value = 24505.385
$ws.Range("L109").Value = 24505.385
$ws.Range("N109").Value = -26585.385
$ws.Range("H113").Value = 911.4286
$ws.Range("I113").Value = 296.6
$ws.Range("J113").Value = 2448.5
$ws.Range("K113").Value = 296.6
$ws.Range("L113").Value = 2448.5
$ws.Range("M113").Value = 1873.4
$ws.Range("N113").Value = -6788.5
$ws.Range("H125").Value = 57662.5
$ws.Range("J125").Value = 57662.5
$ws.Range("L125").Value = 57662.5
$ws.Range("N125").Value = -62582.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 3356.3333
$ws.Range("I51").Value = 3356.3333
$ws.Range("K51").Value = 10068.9999
$ws.Range("M51").Value = -9608.999899999999
$ws.Range("H86").Value = 1179.8
$ws.Range("I86").Value = 1174
$ws.Range("J86").Value = 1181.25
$ws.Range("K86").Value = 3522
$ws.Range("L86").Value = 3543.75
$ws.Range("N86").Value = -5915.75
$ws.Range("M86").Value = -2336
$ws.Range("H89").Value = 1179.8
$ws.Range("I89").Value = 1174
$ws.Range("J89").Value = 1181.25
$ws.Range("K89").Value = 10566
$ws.Range("L89").Value = 10631.25
$ws.Range("N89").Value = -22487.25
$ws.Range("M89").Value = -4638
$ws.Range("H104").Value = 2847.25
$ws.Range("I104").Value = 1700
$ws.Range("J104").Value = 3994.5
$ws.Range("K104").Value = 5100
$ws.Range("L104").Value = 11983.5
$ws.Range("M104").Value = -2479
$ws.Range("N104").Value = -17225.5
$ws.Range("H122").Value = 11112447
$ws.Range("J122").Value = 13890107
$ws.Range("L122").Value = 125010963
$ws.Range("N122").Value = -125015863
$ws.Range("H137").Value = 1955.2667
$ws.Range("I137").Value = 1393.8182
$ws.Range("K137").Value = 4181.4546
$ws.Range("M137").Value = 918.5454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 97728.75
$ws.Range("I70").Value = 142341
$ws.Range("K70").Value = 142341
$ws.Range("M70").Value = -142071
$ws.Range("H73").Value = 97728.75
$ws.Range("I73").Value = 142341
$ws.Range("K73").Value = 142341
$ws.Range("M73").Value = -141405

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 624.5
$ws.Range("I22").Value = 660
$ws.Range("J22").Value = 541.6667
$ws.Range("K22").Value = 660
$ws.Range("L22").Value = 541.6667
$ws.Range("M22").Value = -365
$ws.Range("N22").Value = -1131.6667
$ws.Range("H27").Value = 624.5
$ws.Range("I27").Value = 660
$ws.Range("J27").Value = 541.6667
$ws.Range("K27").Value = 660
$ws.Range("L27").Value = 541.6667
$ws.Range("M27").Value = -553
$ws.Range("N27").Value = -755.6667
$ws.Range("H122").Value = 47623576
$ws.Range("I122").Value = 200002140
$ws.Range("K122").Value = 600006420
$ws.Range("M122").Value = -600003970

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 1000000
$ws.Range("I15").Value = 1000000
$ws.Range("K15").Value = 1000000
$ws.Range("M15").Value = -999712
$ws.Range("H39").Value = 25000
$ws.Range("J39").Value = 25000
$ws.Range("L39").Value = 25000
$ws.Range("N39").Value = -25826
$ws.Range("H40").Value = 25000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 25000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 25000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -25298
$ws.Range("H122").Value = 6271.2856
$ws.Range("I122").Value = 4974.75
$ws.Range("K122").Value = 14924.25
$ws.Range("M122").Value = -12474.25
$ws.Range("H132").Value = 2634.2778
$ws.Range("I132").Value = 2012.8636
$ws.Range("K132").Value = 6038.5908
$ws.Range("M132").Value = -3508.5908
